$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowData = @(
    @{ Row = 2; Cells = @{ "B" = 1.02; "C" = 1.034118942127875; "D" = 1.042667095732968; "E" = 1.037762998858533; "F" = 1.049654476318514; "I" = 1.038031315657763; "J" = 1.039239935007143; "K" = 1.045443136622966; "L" = 1.040552956283314; "M" = 1.052410932474179; "N" = 1.017033415535137 } },
    @{ Row = 3; Cells = @{ "B" = 1.02; "C" = 1.034967773080063; "D" = 1.043350161576656; "E" = 1.038559762122487; "F" = 1.0505145301547; "I" = 1.038225930447641; "J" = 1.039732273808831; "K" = 1.045937415447781; "L" = 1.041159628757937; "M" = 1.053083152467352; "N" = 1.017198539542205 } },
    @{ Row = 4; Cells = @{ "B" = 1.02; "C" = 1.035517665633477; "D" = 1.043792716427896; "E" = 1.039076298368313; "F" = 1.051072079674523; "I" = 1.038350983187015; "J" = 1.040050834643711; "K" = 1.046257132576244; "L" = 1.041552508235686; "M" = 1.053518516804482; "N" = 1.017305331768859 } },
    @{ Row = 5; Cells = @{ "B" = 1.02; "C" = 1.035748992651512; "D" = 1.043978900157349; "E" = 1.039293682349011; "F" = 1.051306719933628; "I" = 1.038403344947149; "J" = 1.040184752745681; "K" = 1.046391512996601; "L" = 1.04171775037355; "M" = 1.053701636496177; "N" = 1.017350213823189 } },
    @{ Row = 6; Cells = @{ "B" = 1.02; "C" = 1.035787842355229; "D" = 1.04401016901294; "E" = 1.039330195636981; "F" = 1.051346131462922; "I" = 1.038412124358576; "J" = 1.040207237863293; "K" = 1.046414074342293; "L" = 1.041745499647319; "M" = 1.053732388478144; "N" = 1.017357748918613 } },
    @{ Row = 7; Cells = @{ "B" = 1.02; "C" = 1.035520756040347; "D" = 1.043795203700942; "E" = 1.039079202153779; "F" = 1.051075213983138; "I" = 1.038351683675043; "J" = 1.040052624084608; "K" = 1.046258928288231; "L" = 1.041554715913945; "M" = 1.053520963297178; "N" = 1.017305931538285 } },
    @{ Row = 8; Cells = @{ "B" = 1.02; "C" = 1.03440567462568; "D" = 1.042897822887037; "E" = 1.038032065337998; "F" = 1.049944919665201; "I" = 1.038097267615288; "J" = 1.039406325647434; "K" = 1.045610203360618; "L" = 1.040757916586355; "M" = 1.052638030038249; "N" = 1.017089230881571 } },
    @{ Row = 9; Cells = @{ "B" = 1.02; "C" = 1.03244574548607; "D" = 1.041320931353256; "E" = 1.036194439151291; "F" = 1.047961222573483; "I" = 1.037642274927668; "J" = 1.038267401740821; "K" = 1.044466248711484; "L" = 1.039356382773918; "M" = 1.051085268730924; "N" = 1.016706982974779 } },
    @{ Row = 10; Cells = @{ "B" = 1.02; "C" = 1.031142571043527; "D" = 1.040272737300312; "E" = 1.034974546255615; "F" = 1.046644265080617; "I" = 1.037334496761056; "J" = 1.037508148713909; "K" = 1.04370313902838; "L" = 1.038423812660725; "M" = 1.050052259110918; "N" = 1.016451913689437 } },
    @{ Row = 11; Cells = @{ "B" = 1.02; "C" = 1.030579118572025; "D" = 1.039819606944622; "E" = 1.034447572937978; "F" = 1.046075337624571; "I" = 1.037200178621017; "J" = 1.037179406007642; "K" = 1.043372608099101; "L" = 1.038020441155865; "M" = 1.049605487738656; "N" = 1.01634141533729 } },
    @{ Row = 12; Cells = @{ "B" = 1.02; "C" = 1.030369953580694; "D" = 1.039651407783764; "E" = 1.034252020955851; "F" = 1.045864213421534; "I" = 1.037150130087729; "J" = 1.037057300383811; "K" = 1.043249820563714; "L" = 1.03787067820606; "M" = 1.049439617969828; "N" = 1.016300364029039 } },
    @{ Row = 13; Cells = @{ "B" = 1.02; "C" = 1.030414814439139; "D" = 1.039687481904488; "E" = 1.03429395892115; "F" = 1.04590949116335; "I" = 1.037160872752417; "J" = 1.03708349225154; "K" = 1.043276159502794; "L" = 1.03790279980342; "M" = 1.049475193906769; "N" = 1.016309169996946 } },
    @{ Row = 14; Cells = @{ "B" = 1.02; "C" = 1.030561826342244; "D" = 1.039805701221327; "E" = 1.034431404676484; "F" = 1.046057881925156; "I" = 1.037196044788288; "J" = 1.037169312628911; "K" = 1.043362458713661; "L" = 1.03800806032469; "M" = 1.049591775229993; "N" = 1.016338022169595 } },
    @{ Row = 15; Cells = @{ "B" = 1.02; "C" = 1.030652422015943; "D" = 1.039878555153814; "E" = 1.034516114690638; "F" = 1.046149337033772; "I" = 1.037217694683831; "J" = 1.037222189969627; "K" = 1.043415628738398; "L" = 1.038072923752418; "M" = 1.049663615613376; "N" = 1.016355797991867 } },
    @{ Row = 16; Cells = @{ "B" = 1.02; "C" = 1.031179983161396; "D" = 1.040302825907776; "E" = 1.035009546195955; "F" = 1.046682050975472; "I" = 1.03734338898797; "J" = 1.037529966786237; "K" = 1.043725073294478; "L" = 1.038450592441141; "M" = 1.050081921128703; "N" = 1.016459246052019 } },
    @{ Row = 17; Cells = @{ "B" = 1.02; "C" = 1.03151113181228; "D" = 1.040569160334534; "E" = 1.035319398248218; "F" = 1.047016564121236; "I" = 1.037421953559849; "J" = 1.037723032963719; "K" = 1.043919153995877; "L" = 1.038687612193088; "M" = 1.050344455775576; "N" = 1.016524122780375 } },
    @{ Row = 18; Cells = @{ "B" = 1.02; "C" = 1.031704365260607; "D" = 1.040724580423223; "E" = 1.035500250034234; "F" = 1.047211807582789; "I" = 1.037467677767238; "J" = 1.0378356469369; "K" = 1.044032348307292; "L" = 1.03882590391982; "M" = 1.050497638689632; "N" = 1.016561959262155 } },
    @{ Row = 19; Cells = @{ "B" = 1.02; "C" = 1.031770266402547; "D" = 1.040777586768959; "E" = 1.035561936173676; "F" = 1.047278402143695; "I" = 1.037483251358136; "J" = 1.037874045649278; "K" = 1.044070942961091; "L" = 1.038873064924069; "M" = 1.050549878689146; "N" = 1.016574859659751 } },
    @{ Row = 20; Cells = @{ "B" = 1.02; "C" = 1.031475594405627; "D" = 1.040540577739677; "E" = 1.035286141596384; "F" = 1.04698066081423; "I" = 1.037413534785252; "J" = 1.037702318601431; "K" = 1.043898331965034; "L" = 1.038662177850124; "M" = 1.05031628301663; "N" = 1.016517162636927 } },
    @{ Row = 21; Cells = @{ "B" = 1.02; "C" = 1.030518531500792; "D" = 1.03977088542819; "E" = 1.034390925085086; "F" = 1.046014178986457; "I" = 1.037185691823225; "J" = 1.037144040545062; "K" = 1.043337046112841; "L" = 1.037977061858329; "M" = 1.049557442686905; "N" = 1.016329526121684 } },
    @{ Row = 22; Cells = @{ "B" = 1.02; "C" = 1.029917519824822; "D" = 1.03928760761669; "E" = 1.033829163479711; "F" = 1.045407676124425; "I" = 1.037041531076968; "J" = 1.036793053036208; "K" = 1.042984065198702; "L" = 1.03754669168814; "M" = 1.049080799149503; "N" = 1.016211509588574 } },
    @{ Row = 23; Cells = @{ "B" = 1.02; "C" = 1.03023605760594; "D" = 1.039543739341332; "E" = 1.034126859387888; "F" = 1.045729083930684; "I" = 1.037118039146112; "J" = 1.036979115466887; "K" = 1.043171194039256; "L" = 1.0377748016567; "M" = 1.049333431871842; "N" = 1.016274076213344 } },
    @{ Row = 24; Cells = @{ "B" = 1.02; "C" = 1.031491651982865; "D" = 1.040553492762581; "E" = 1.035301168474275; "F" = 1.046996883579192; "I" = 1.037417339179815; "J" = 1.037711678524889; "K" = 1.043907740574563; "L" = 1.038673670404952; "M" = 1.050329012916652; "N" = 1.016520307641712 } },
    @{ Row = 25; Cells = @{ "B" = 1.02; "C" = 1.032951833708251; "D" = 1.04172806259155; "E" = 1.036668602378271; "F" = 1.048473093268761; "I" = 1.037760688740197; "J" = 1.038561840951789; "K" = 1.044762077465112; "L" = 1.039718404391269; "M" = 1.051486320099594; "N" = 1.016805847055712 } }
)

foreach ($r in $rowData) {
    $rowNum = $r.Row
    foreach ($col in $r.Cells.Keys) {
        $ws.Range("$col$rowNum").Value = $r.Cells[$col]
    }
}

Write-Output "Updated $($rowData.Count) rows"
